$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 currently holds "192.168.1.113" -> becomes "127.0.0.1"
$ws.Range("C2").Value = "127.0.0.1"

# E2 currently holds "192.168.0.24" and picks up C2's old number format/style (s="1")
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("E2").NumberFormat = "@"

# Update the active selection to E2 (was C2)
$ws.Range("E2").Select()
